# "call not support '$' replace '.'"
#
# The Copy1/Copy2 ellipse labels and the SystemA ellipse labels on the
# "Copy" system diagrams use '$' as a separator between the (copied)
# system name and the function name, e.g. "Copy1 $Func1" / "SystemA$Func2".
# Because the runtime that consumes these labels does not support '$' as
# a call-separator, every occurrence of '$' in those labels is replaced
# with '.'.

$p = $ppt.ActivePresentation

# --- Slide 2 ("S101"): the generic "Copy1 / Copy2" template shapes -------
$slide2 = $p.Slides.Item(2)
$grp7_s2  = $slide2.Shapes.Item("그룹 7")
$grp2_s2  = $grp7_s2.GroupItems.Item("그룹 2")

$copy1 = $grp2_s2.GroupItems.Item("타원 10")
$copy1.TextFrame.TextRange.Text = "Copy1 .Func1"

$copy2 = $grp2_s2.GroupItems.Item("타원 11")
$copy2.TextFrame.TextRange.Text = "Copy2 .Func2"

# --- Slide 3 ("S102"): the concrete "SystemA" instance shapes ------------
$slide3 = $p.Slides.Item(3)
$grp7_s3  = $slide3.Shapes.Item("그룹 7")
$grp2_s3  = $grp7_s3.GroupItems.Item("그룹 2")

# "타원 10" originally holds the text split across two runs
# ("SystemA $" + "Func1"). Collapse it to a short placeholder first so the
# whole label becomes a single run, then assign the final text - this
# produces one clean run instead of leaving stray run fragments behind.
$systemAFunc1 = $grp2_s3.GroupItems.Item("타원 10")
$systemAFunc1.TextFrame.TextRange.Text = "x"
$systemAFunc1.TextFrame.TextRange.Text = "SystemA .Func1"

# "타원 11" already holds its text ("SystemA$Func2") in a single run.
$systemAFunc2 = $grp2_s3.GroupItems.Item("타원 11")
$systemAFunc2.TextFrame.TextRange.Text = "SystemA.Func2"
